$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Include from International Cl" sheet to "Include #0" ---
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Name = "Include #0"

# --- 2. Metadata sheet updates ---
$ws1 = $wb.Worksheets.Item(1)

# Give new row 15 the same cell formatting as row 14 (format-only copy reuses
# the existing style index instead of minting a new one)
$ws1.Range("A14:B14").Copy()
$ws1.Range("A15:B15").PasteSpecial(-4122)

# Capture current (pre-shift) values of rows 11-14 (both columns) before
# overwriting anything
$a11 = $ws1.Range("A11").Value2
$b11 = $ws1.Range("B11").Value2
$a12 = $ws1.Range("A12").Value2
$b12 = $ws1.Range("B12").Value2
$a13 = $ws1.Range("A13").Value2
$b13 = $ws1.Range("B13").Value2
$a14 = $ws1.Range("A14").Value2
$b14 = $ws1.Range("B14").Value2

# Shift rows 11-14 down to rows 12-15
$ws1.Range("A15").Value = $a14
$ws1.Range("B15").Value = $b14

$ws1.Range("A14").Value = $a13
$ws1.Range("B14").Value = $b13

$ws1.Range("A13").Value = $a12
$ws1.Range("B13").Value = $b12

$ws1.Range("A12").Value = $a11
$ws1.Range("B12").Value = $b11

# New row 11: "Jurisdiction" property with an (empty) value
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

# --- 3. Update the Date value (row 8, column B) ---
$ws1.Range("B8").Value = "2024-09-17T19:55:11+00:00"
